$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
if (-not $ws) { $ws = $wb.ActiveSheet }

# The CasesTab query (row 2, column B) is being edited: the trailing
# "Cohort" column is dropped from the returned result set, so the last
# RETURN line (coalesce(co.cohort_description, '') AS `Cohort`) is removed.
$newCasesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`nMATCH (c)<--(diag:diagnosis)`nMATCH (samp:sample)-->(c) `n WHERE samp.specific_sample_pathology IN [""Undefined""]  `nOPTIONAL MATCH (co:cohort)<-[*]-(c)`n  WITH DISTINCT c, s, demo, diag, co`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $newCasesQuery

# Update the selection/active cell shown when the sheet is reopened.
$ws.Activate()
$ws.Range("B2").Select()

$wb.Save()
